# Adds new-country figures (China / Brazil / Mexico / South Korea -> columns N:Q)
# to the "Specificities" sheet. Columns N1:Q1 (country names) already existed;
# this fills in the corresponding data rows that were still empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specificities")

# ---------------------------------------------------------------------------
# Row 7: currency symbol header (¥ / R$ / $ / ₩)
# ---------------------------------------------------------------------------
$ws.Range("N7").Value = "¥"
$ws.Range("O7").Value = "R$"
$ws.Range("P7").Value = "$"
$ws.Range("Q7").Value = "₩"

# ---------------------------------------------------------------------------
# Row 8: currency conversion rate (currency in $ as of 5/5/21)
# ---------------------------------------------------------------------------
$ws.Range("N8").Value = 6.37
$ws.Range("O8").Value = 5.23
$ws.Range("P8").Value = 19.920000000000002
$ws.Range("Q8").Value = 1113.8499999999999

# ---------------------------------------------------------------------------
# Rows 9-12: job creation/destruction (rounded figures, text)
# ---------------------------------------------------------------------------
$ws.Range("N9").Value = "12,480k"
$ws.Range("O9").Value = "1,445k"
$ws.Range("P9").Value = "815k"
$ws.Range("Q9").Value = "1,110k"

$ws.Range("N10").Value = "3,375k"
$ws.Range("O10").Value = "825k"
$ws.Range("P10").Value = "550k"
$ws.Range("Q10").Value = "205k"

$ws.Range("N11").Value = "12.5M"
$ws.Range("O11").Value = "1.5M"
$ws.Range("P11").Value = "800k"
$ws.Range("Q11").Value = "1M"

$ws.Range("N12").Value = "3.5M"
$ws.Range("O12").Value = "800k"
$ws.Range("P12").Value = "550k"
$ws.Range("Q12").Value = "200k"

# ---------------------------------------------------------------------------
# Row 51: carbon price ($/tCO2)
# ---------------------------------------------------------------------------
$ws.Range("N51").Value = 45
$ws.Range("O51").Value = 45
$ws.Range("P51").Value = 45
$ws.Range("Q51").Value = 45

# ---------------------------------------------------------------------------
# Row 52: fossil CO2 emissions (2017, MtCO2) - copy M52's numeric format first
# ---------------------------------------------------------------------------
$ws.Range("M52").Copy()
$ws.Range("N52:Q52").PasteSpecial(-4122)
$ws.Range("N52").Value = 10877
$ws.Range("O52").Value = 493
$ws.Range("P52").Value = 507
$ws.Range("Q52").Value = 673

# ---------------------------------------------------------------------------
# Row 53: adult population (2020)
# ---------------------------------------------------------------------------
$ws.Range("N53").Value = 1128677232
$ws.Range("O53").Value = 159837762.59999999
$ws.Range("P53").Value = 92799575.799999997
$ws.Range("Q53").Value = 43190072.399999999

# ---------------------------------------------------------------------------
# Row 46: 20.7 $ global tax pc = $B$46 * <country column 8>
# ---------------------------------------------------------------------------
$ws.Range("M46").Copy()
$ws.Range("N46:Q46").PasteSpecial(-4122)
$ws.Range("N46").Formula = '=$B$46*N8'
$ws.Range("O46").Formula = '=$B$46*O8'
$ws.Range("P46").Formula = '=$B$46*P8'
$ws.Range("Q46").Formula = '=$B$46*Q8'

# ---------------------------------------------------------------------------
# Row 54: gasoline price increase in $/liter = 2.5*<row51>/1000
# ---------------------------------------------------------------------------
$ws.Range("M54").Copy()
$ws.Range("N54:Q54").PasteSpecial(-4122)
$ws.Range("N54").Formula = '=2.5*N51/1000'
$ws.Range("O54").Formula = '=2.5*O51/1000'
$ws.Range("P54").Formula = '=2.5*P51/1000'
$ws.Range("Q54").Formula = '=2.5*Q51/1000'

# ---------------------------------------------------------------------------
# Row 55: gasoline price increase in LCU/liter = <row54> * <row8>
# ---------------------------------------------------------------------------
$ws.Range("M55").Copy()
$ws.Range("N55:Q55").PasteSpecial(-4122)
$ws.Range("N55").Formula = '=N54*N8'
$ws.Range("O55").Formula = '=O54*O8'
$ws.Range("P55").Formula = '=P54*P8'
$ws.Range("Q55").Formula = '=Q54*Q8'

# ---------------------------------------------------------------------------
# Row 56: transfer per adult in LCU/year
# ---------------------------------------------------------------------------
$ws.Range("M56").Copy()
$ws.Range("N56:Q56").PasteSpecial(-4122)
$ws.Range("N56").Formula = '=0.8*0.8*N52*N51*1000000*N8/N53'
$ws.Range("O56").Formula = '=0.8*0.8*O52*O51*1000000*O8/O53'
$ws.Range("P56").Formula = '=0.8*0.8*P52*P51*1000000*P8/P53'
$ws.Range("Q56").Formula = '=0.8*0.8*Q52*Q51*1000000*Q8/Q53'
